# Implement heterogenous mixing using socialising params
#
# This edit replaces the old "mixing_factor_cc" / "mixing_factor_ca" parameter
# rows on the "constant" sheet with new "child_socialising" / "elderly_socialising"
# parameter rows (including new distribution bounds), and widens the uniform
# prior on raw_transmission_rate's distribution parameters.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constant")

# --- Row 2: raw_transmission_rate - widen distribution bounds ---
$ws.Range("D2").Value = 0.001
$ws.Range("E2").Value = 0.01

# --- Rows 3 & 4: replace mixing_factor_cc/ca with the new socialising params ---
# Set the string-valued cells first, and in this particular order, so that the
# new entries land in the shared-string table in the same order as the source
# workbook (parameter names first, then the "65+" description, then the "0-14"
# description).
$ws.Range("A3").Value = "child_socialising"
$ws.Range("A4").Value = "elderly_socialising"
$ws.Range("G4").Value = "Social activity level for 65+ years old (relative to 15-64 years-old)"
$ws.Range("G3").Value = "Social activity level for 0-14 years old (relative to 15-64 years-old)"

# Row 3 - child_socialising
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "uniform"
$ws.Range("D3").Value = 0.2
$ws.Range("E3").Value = 1
$ws.Range("H3").Value = "m_{cc}"

# Row 4 - elderly_socialising
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "uniform"
$ws.Range("D4").Value = 0.2
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = "m_{ca}"

# --- View state: constant sheet selection moves to D4 ---
$ws.Activate()
$ws.Range("D4").Select()
